# Applies the diff: insert two new data rows (407 and 408) into the sheet,
# pushing the existing rows 407..443 down to 409..445, and fill the new
# rows with their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 407.
$ws.Rows("407:408").Insert()

# Row 407: Copenhague / Primera
$ws.Cells.Item(407, 1).Value = 7
$ws.Cells.Item(407, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(407, 3).Value = "Ñuble"
$ws.Cells.Item(407, 4).Value = 45106
$ws.Cells.Item(407, 5).Value = 16
$ws.Cells.Item(407, 6).Value = 100112006
$ws.Cells.Item(407, 7).Value = "Repollo"
$ws.Cells.Item(407, 8).Value = "Copenhague"
$ws.Cells.Item(407, 9).Value = "Primera"
$ws.Cells.Item(407, 10).Value = 100
$ws.Cells.Item(407, 11).Value = 1000
$ws.Cells.Item(407, 12).Value = 1000
$ws.Cells.Item(407, 13).Value = 1000
$ws.Cells.Item(407, 14).Value = "$/unidad"
$ws.Cells.Item(407, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(407, 16).Value = 1000
$ws.Cells.Item(407, 17).Value = 1
$ws.Cells.Item(407, 18).Value = "Hortaliza"

# Row 408: Crespo record / Primera
$ws.Cells.Item(408, 1).Value = 7
$ws.Cells.Item(408, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(408, 3).Value = "Ñuble"
$ws.Cells.Item(408, 4).Value = 45106
$ws.Cells.Item(408, 5).Value = 16
$ws.Cells.Item(408, 6).Value = 100112006
$ws.Cells.Item(408, 7).Value = "Repollo"
$ws.Cells.Item(408, 8).Value = "Crespo record"
$ws.Cells.Item(408, 9).Value = "Primera"
$ws.Cells.Item(408, 10).Value = 150
$ws.Cells.Item(408, 11).Value = 1000
$ws.Cells.Item(408, 12).Value = 1000
$ws.Cells.Item(408, 13).Value = 1000
$ws.Cells.Item(408, 14).Value = "$/unidad"
$ws.Cells.Item(408, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(408, 16).Value = 1000
$ws.Cells.Item(408, 17).Value = 1
$ws.Cells.Item(408, 18).Value = "Hortaliza"

# Make sure the date column keeps the expected date format (it should
# already inherit this from the surrounding rows after the row insert).
$ws.Range("D407:D408").NumberFormat = "YYYY-MM-DD HH:MM:SS"
